# "Adding pos. mode adducts for PUAs"
#
# 1. Record a new abundance rank (1) for the [M+H]+ adduct of PUA (row 14,
#    column M) on the "Adduct ion hierarchies" sheet.
# 2. Log the change on the "Notes" sheet: a new dated row describing the
#    edit, authored by "JRC", formatted like the existing history rows.
# 3. Leave the selection/active sheet the way the author left the workbook:
#    cell M15 selected on "Adduct ion hierarchies", then the "Notes" sheet
#    activated with B16 selected.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Adduct ion hierarchies")
$ws2 = $wb.Worksheets.Item("Notes")

# --- 1. New adduct-ion-hierarchy data point ---------------------------------
# PUA row (row 14) gets a rank of 1 for the [M+H]+ adduct (column M).
$ws1.Range("M14").Value = 1

# --- 2. Append a new file-history entry on the Notes sheet -------------------
$ws2.Range("A15").Value = 42403   # 2016-02-03, serial date like the rows above
$ws2.Range("A14").Copy() | Out-Null
$ws2.Range("A15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats: reuse the date number format
$ws2.Range("B15").Value = "Added positive ion mode adduct ([M+H]+) for PUAs"
$ws2.Range("C15").Value = "JRC"

# --- 3. Restore selection / active sheet -------------------------------------
$ws1.Range("M15").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("B16").Select() | Out-Null
